$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 12 and 13 (data now only spans to row 11; corresponding
# "Inflammatory-Mac" target-cluster rows are removed)
$ws.Rows("12:13").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.500469666666667
$ws.Range("H2").Value = 4.501409000000001
$ws.Range("I2").Value = 0.9943843705197677
$ws.Range("J2").Value = 0.9943843705197678
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2328916666666666
$ws.Range("N2").Value = 0.6986749999999999
$ws.Range("O2").Value = 0.01421300418632399
$ws.Range("P2").Value = 0.01421300418632399
$ws.Range("Q2").Value = 0.3494468814527777
$ws.Range("R2").Value = 3.145021933075
$ws.Range("S2").Value = 0.01413318922101261
$ws.Range("T2").Value = 0.01413318922101261

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.500469666666667
$ws.Range("H3").Value = 4.501409000000001
$ws.Range("I3").Value = 0.9943843705197677
$ws.Range("J3").Value = 0.9943843705197678
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.57672866666667
$ws.Range("N3").Value = 43.730186
$ws.Range("O3").Value = 0.8895943273864486
$ws.Range("P3").Value = 0.8895943273864487
$ws.Range("Q3").Value = 21.87193920356378
$ws.Range("R3").Value = 196.8474528320741
$ws.Range("S3").Value = 0.8845986952561299
$ws.Range("T3").Value = 0.88459869525613

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.500469666666667
$ws.Range("H4").Value = 4.501409000000001
$ws.Range("I4").Value = 0.9943843705197677
$ws.Range("J4").Value = 0.9943843705197678
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.423796666666667
$ws.Range("N4").Value = 4.27139
$ws.Range("O4").Value = 0.08689202268783405
$ws.Range("P4").Value = 0.08689202268783405
$ws.Range("Q4").Value = 2.136363709834445
$ws.Range("R4").Value = 19.22727338851
$ws.Range("S4").Value = 0.08640406928363124
$ws.Range("T4").Value = 0.08640406928363124

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.500469666666667
$ws.Range("H5").Value = 4.501409000000001
$ws.Range("I5").Value = 0.9943843705197677
$ws.Range("J5").Value = 0.9943843705197678
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.136948
$ws.Range("N5").Value = 0.410844
$ws.Range("O5").Value = 0.008357716380185487
$ws.Range("P5").Value = 0.008357716380185487
$ws.Range("Q5").Value = 0.2054863199106667
$ws.Range("R5").Value = 1.849376879196
$ws.Range("S5").Value = 0.008310782541693497
$ws.Range("T5").Value = 0.008310782541693498

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.500469666666667
$ws.Range("H6").Value = 4.501409000000001
$ws.Range("I6").Value = 0.9943843705197677
$ws.Range("J6").Value = 0.9943843705197678
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.01545066666666667
$ws.Range("N6").Value = 0.046352
$ws.Range("O6").Value = 0.0009429293592077716
$ws.Range("P6").Value = 0.0009429293592077717
$ws.Range("Q6").Value = 0.02318325666311111
$ws.Range("R6").Value = 0.208649309968
$ws.Range("S6").Value = 0.0009376342173004279
$ws.Range("T6").Value = 0.0009376342173004281

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.008473666666666666
$ws.Range("H7").Value = 0.025421
$ws.Range("I7").Value = 0.005615629480232302
$ws.Range("J7").Value = 0.005615629480232303
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2328916666666666
$ws.Range("N7").Value = 0.6986749999999999
$ws.Range("O7").Value = 0.01421300418632399
$ws.Range("P7").Value = 0.01421300418632399
$ws.Range("Q7").Value = 0.001973446352777777
$ws.Range("R7").Value = 0.017761017175
$ws.Range("S7").Value = [double]"7.981496531138614E-05"
$ws.Range("T7").Value = [double]"7.981496531138616E-05"

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008473666666666666
$ws.Range("H8").Value = 0.025421
$ws.Range("I8").Value = 0.005615629480232302
$ws.Range("J8").Value = 0.005615629480232303
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.57672866666667
$ws.Range("N8").Value = 43.730186
$ws.Range("O8").Value = 0.8895943273864486
$ws.Range("P8").Value = 0.8895943273864487
$ws.Range("Q8").Value = 0.1235183398117778
$ws.Range("R8").Value = 1.111665058306
$ws.Range("S8").Value = 0.004995632130318767
$ws.Range("T8").Value = 0.004995632130318769

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008473666666666666
$ws.Range("H9").Value = 0.025421
$ws.Range("I9").Value = 0.005615629480232302
$ws.Range("J9").Value = 0.005615629480232303
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.423796666666667
$ws.Range("N9").Value = 4.27139
$ws.Range("O9").Value = 0.08689202268783405
$ws.Range("P9").Value = 0.08689202268783405
$ws.Range("Q9").Value = 0.01206477835444444
$ws.Range("R9").Value = 0.10858300519
$ws.Range("S9").Value = 0.0004879534042028149
$ws.Range("T9").Value = 0.0004879534042028151

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008473666666666666
$ws.Range("H10").Value = 0.025421
$ws.Range("I10").Value = 0.005615629480232302
$ws.Range("J10").Value = 0.005615629480232303
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.136948
$ws.Range("N10").Value = 0.410844
$ws.Range("O10").Value = 0.008357716380185487
$ws.Range("P10").Value = 0.008357716380185487
$ws.Range("Q10").Value = 0.001160451702666666
$ws.Range("R10").Value = 0.010444065324
$ws.Range("S10").Value = [double]"4.693383849199003E-05"
$ws.Range("T10").Value = [double]"4.693383849199003E-05"

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.008473666666666666
$ws.Range("H11").Value = 0.025421
$ws.Range("I11").Value = 0.005615629480232302
$ws.Range("J11").Value = 0.005615629480232303
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.01545066666666667
$ws.Range("N11").Value = 0.046352
$ws.Range("O11").Value = 0.0009429293592077716
$ws.Range("P11").Value = 0.0009429293592077717
$ws.Range("Q11").Value = 0.0001309237991111111
$ws.Range("R11").Value = 0.001178314192
$ws.Range("S11").Value = [double]"5.295141907343717E-06"
$ws.Range("T11").Value = [double]"5.295141907343718E-06"
